$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.79"
$ws.Range("E2").Value = "'-4.46%"

$ws.Range("D3").Value = "'35.29"
$ws.Range("E3").Value = "'-2.22%"

$ws.Range("D4").Value = "'5.086"
$ws.Range("E4").Value = "'-1.17%"

$ws.Range("D5").Value = "'0.08003"
$ws.Range("E5").Value = "'-3.08%"

$ws.Range("D6").Value = "'1.936"
$ws.Range("E6").Value = "'-9.96%"

$ws.Range("D7").Value = "'4.062"
$ws.Range("E7").Value = "'-1.91%"

$ws.Range("D8").Value = "'7.771"
$ws.Range("E8").Value = "'-3.17%"

$ws.Range("D9").Value = "'2.958"
$ws.Range("E9").Value = "'5.74%"

$ws.Range("D10").Value = "'0.9219"
$ws.Range("E10").Value = "'-0.60%"

$ws.Range("D11").Value = "'0.1223"
$ws.Range("E11").Value = "'19.44%"

$ws.Range("D12").Value = "'0.1854"
$ws.Range("E12").Value = "'-1.67%"

$ws.Range("D13").Value = "'0.09738"
$ws.Range("E13").Value = "'4.80%"

$ws.Range("D14").Value = "'0.03637"
$ws.Range("E14").Value = "'0.64%"

$ws.Range("D15").Value = "'0.09864"
$ws.Range("E15").Value = "'-0.61%"

$ws.Range("D16").Value = "'0.001397"

$ws.Range("D17").Value = "'0.005793"
$ws.Range("E17").Value = "'1.98%"

$ws.Range("D18").Value = "'3.498"
$ws.Range("E18").Value = "'1.21%"

$ws.Range("E19").Value = "'0.87%"

$ws.Range("D20").Value = "'0.1302"
$ws.Range("E20").Value = "'-1.41%"

$ws.Range("D21").Value = "'5.039"
$ws.Range("E21").Value = "'-2.91%"

$ws.Range("D22").Value = "'0.2467"
$ws.Range("E22").Value = "'12.54%"

$ws.Range("D23").Value = "'0.04528"
$ws.Range("E23").Value = "'-1.43%"

$ws.Range("E24").Value = "'-2.53%"

$ws.Range("D25").Value = "'0.004838"
$ws.Range("E25").Value = "'2.22%"

$ws.Range("D26").Value = "'0.0001252"
$ws.Range("E26").Value = "'0.02%"

$ws.Range("E27").Value = "'-6.93%"

$ws.Range("D39").Value = "'0.01930"
$ws.Range("E39").Value = "'-3.40%"

$ws.Range("D40").Value = "'0.04720"
$ws.Range("E40").Value = "'-4.59%"

$ws.Range("D41").Value = "'0.007674"
$ws.Range("E41").Value = "'-1.54%"

$ws.Range("D42").Value = "'0.009726"
$ws.Range("E42").Value = "'24.12%"

$ws.Range("E43").Value = "'-5.16%"

$ws.Range("D44").Value = "'0.002113"
$ws.Range("E44").Value = "'-1.38%"

$ws.Range("D45").Value = "'0.01013"
$ws.Range("E45").Value = "'-13.57%"

$ws.Range("D46").Value = "'0.00006282"
$ws.Range("E46").Value = "'-2.90%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.01%"

$ws.Range("E48").Value = "'91.50%"

$ws.Range("D49").Value = "'0.001490"
$ws.Range("E49").Value = "'-21.67%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.01%"

$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.01%"
